$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 4 (the row currently holding 005206566 / LEVI),
# shifting that row and everything below it down by one.
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row with the account data.
# The account number column holds text (with a leading zero), so momentarily
# force the cell to text format before writing the value - otherwise Excel
# would coerce the numeric-looking string and drop the leading zero. The
# number format is cleared again afterwards so the cell is left with the
# same (default) styling as its neighbours, only its text content differs.
$accountCell = $ws.Cells.Item(4, 1)
$accountCell.NumberFormat = "@"
$accountCell.Value = "002606448"
$accountCell.ClearFormats()

$ws.Cells.Item(4, 2).Value = "MARCUS"
$ws.Cells.Item(4, 3).Value = 58000
